$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = -3.5055773183210661
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2.2934989013837068
$ws.Range("E2").Value = 2.7364738369104487

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 7.1887756209839324
$ws.Range("D3").Value = 8.4695263851135536
$ws.Range("E3").Value = 1.0364623561658632

# Update the selected range to match the new reduced data extent
$ws.Range("B1:E3").Select()
